$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-text number format on price cells whose new values would
# otherwise be auto-parsed by Excel as numbers (losing formatting such as
# trailing zeros or switching to scientific notation).
$textCells = @("D5", "D7", "D9", "D10", "D11", "D13", "D14", "D15", "D16", "D18", "D20", "D22", "D25", "D26", "D27", "D28", "D30", "D31", "D32", "D33", "D34", "D35", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "26.862.94"
$ws.Range("E2").Value = "  -1.78%  "

# Row 3
$ws.Range("D3").Value = "1.810.32"
$ws.Range("E3").Value = "  -0.89%  "

# Row 4
$ws.Range("E4").Value = "  +0.11%  "

# Row 5
$ws.Range("D5").Value = "309.75"
$ws.Range("E5").Value = "  -1.43%  "

# Row 6
$ws.Range("E6").Value = "  +0.11%  "

# Row 7
$ws.Range("D7").Value = "0.4649"
$ws.Range("E7").Value = "  -1.05%  "

# Row 8
$ws.Range("E8").Value = "  -2.33%  "

# Row 9
$ws.Range("D9").Value = "0.07371"
$ws.Range("E9").Value = "  -0.68%  "

# Row 10
$ws.Range("D10").Value = "0.8758"
$ws.Range("E10").Value = "  -0.24%  "

# Row 11
$ws.Range("D11").Value = "20.46"
$ws.Range("E11").Value = "  -1.66%  "

# Row 12
$ws.Range("D12").Value = "1.851.89"
$ws.Range("E12").Value = "  +1.37%  "

# Row 13
$ws.Range("D13").Value = "5.364"
$ws.Range("E13").Value = "  -1.34%  "

# Row 14
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "6.509"
$ws.Range("E14").Value = "  -2.95%  "

# Row 15
$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D15").Value = "0.07059"
$ws.Range("E15").Value = "  -0.40%  "

# Row 16
$ws.Range("D16").Value = "91.60"
$ws.Range("E16").Value = "  -1.82%  "

# Row 17
$ws.Range("E17").Value = "  +0.13%  "

# Row 18
$ws.Range("D18").Value = "0.000008704"
$ws.Range("E18").Value = "  -1.22%  "

# Row 19
$ws.Range("E19").Value = "  +0.10%  "

# Row 20
$ws.Range("D20").Value = "14.74"
$ws.Range("E20").Value = "  -1.89%  "

# Row 21
$ws.Range("D21").Value = "26.862.97"
$ws.Range("E21").Value = "  -1.78%  "

# Row 22
$ws.Range("D22").Value = "5.316"
$ws.Range("E22").Value = "  -0.39%  "

# Row 23
$ws.Range("E23").Value = "  -3.79%  "

# Row 24
$ws.Range("D24").Value = "2.008.96"
$ws.Range("E24").Value = "  -1.99%  "

# Row 25
$ws.Range("D25").Value = "1.901"
$ws.Range("E25").Value = "  -2.09%  "

# Row 26
$ws.Range("D26").Value = "151.61"
$ws.Range("E26").Value = "  +0.28%  "

# Row 27
$ws.Range("D27").Value = "18.38"
$ws.Range("E27").Value = "  -1.26%  "

# Row 28
$ws.Range("D28").Value = "2.150"
$ws.Range("E28").Value = "  -4.69%  "

# Row 29
$ws.Range("E29").Value = "  -0.56%  "

# Row 30
$ws.Range("D30").Value = "115.96"
$ws.Range("E30").Value = "  -1.09%  "

# Row 31
$ws.Range("D31").Value = "0.08910"
$ws.Range("E31").Value = "  -0.60%  "

# Row 32
$ws.Range("D32").Value = "0.7540"
$ws.Range("E32").Value = "  -5.12%  "

# Row 33
$ws.Range("D33").Value = "1.153"
$ws.Range("E33").Value = "  -3.52%  "

# Row 34
$ws.Range("D34").Value = "4.463"
$ws.Range("E34").Value = "  -1.94%  "

# Row 35
$ws.Range("D35").Value = "2.919"
$ws.Range("E35").Value = "  -0.51%  "

# Row 36
$ws.Range("E36").Value = "  +0.17%  "

# Row 37
$ws.Range("E37").Value = "  -0.07%  "

# Row 38
$ws.Range("E38").Value = "  -0.72%  "

# Row 39
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "0.05259"
$ws.Range("E39").Value = "  +0.11%  "

# Row 40
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "2.439"
$ws.Range("E40").Value = "  +2.19%  "

# Row 41
$ws.Range("D41").Value = "2.918"
$ws.Range("E41").Value = "  +0.87%  "

# Row 42
$ws.Range("D42").Value = "0.5302"
$ws.Range("E42").Value = "  -0.69%  "

# Row 43
$ws.Range("D43").Value = "7.169"
$ws.Range("E43").Value = "  -2.04%  "

# Row 44
$ws.Range("D44").Value = "0.1661"
$ws.Range("E44").Value = "  -2.62%  "

# Row 45
$ws.Range("D45").Value = "8.453"
$ws.Range("E45").Value = "  -2.37%  "

# Row 46
$ws.Range("D46").Value = "0.4943"
$ws.Range("E46").Value = "  -3.41%  "

# Row 47
$ws.Range("D47").Value = "10.35"
$ws.Range("E47").Value = "  -2.49%  "

# Row 48
$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D48").Value = "1.001"
$ws.Range("E48").Value = "  +0.20%  "

# Row 49
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "1.673"
$ws.Range("E49").Value = "  -0.65%  "

# Row 50
$ws.Range("D50").Value = "103.06"
$ws.Range("E50").Value = "  -2.31%  "

# Row 51
$ws.Range("D51").Value = "0.06286"
$ws.Range("E51").Value = "  -1.60%  "
